$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Jaden Ivey, PG,SG, Detroit Pistons -> Payton Pritchard, PG, Boston Celtics
$ws.Range("A2").Value = "Payton Pritchard"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Boston Celtics"

# Row 6: Anthony Edwards, SG,SF, Minnesota Timberwolves -> Jayson Tatum, SF,PF, Boston Celtics
$ws.Range("A6").Value = "Jayson Tatum"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Boston Celtics"

# Row 7: Jaden McDaniels, SF,PF, Minnesota Timberwolves -> Jaren Jackson Jr., PF,C, Memphis Grizzlies
$ws.Range("A7").Value = "Jaren Jackson Jr."
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Memphis Grizzlies"

# Row 11: Wendell Carter Jr., C, Orlando Magic -> Draymond Green, PF,C, Golden State Warriors
$ws.Range("A11").Value = "Draymond Green"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Golden State Warriors"

# Row 14: Jaren Jackson Jr., PF,C, Memphis Grizzlies -> Jaden Ivey, PG,SG, Detroit Pistons
$ws.Range("A14").Value = "Jaden Ivey"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Detroit Pistons"

# Row 15: Draymond Green, PF,C, Golden State Warriors -> Anthony Edwards, SG,SF, Minnesota Timberwolves
$ws.Range("A15").Value = "Anthony Edwards"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Minnesota Timberwolves"

# Row 16: Jayson Tatum, SF,PF, Boston Celtics -> Wendell Carter Jr., C, Orlando Magic
$ws.Range("A16").Value = "Wendell Carter Jr."
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Orlando Magic"
